$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 203, shifting existing rows 203:287 down to 204:288.
$ws.Rows.Item(203).Insert()

# Populate the newly inserted row 203 with the new record.
$ws.Range("A203").Value = 10
$ws.Range("B203").Value = "Vega Modelo de Temuco"
$ws.Range("C203").Value = "La Araucanía"
$ws.Range("D203").Value = 44636
$ws.Range("E203").Value = 9
$ws.Range("F203").Value = 100112044
$ws.Range("G203").Value = "Perejil"
$ws.Range("H203").Value = "Sin especificar"
$ws.Range("I203").Value = "Primera"
$ws.Range("J203").Value = 20
$ws.Range("K203").Value = 5000
$ws.Range("L203").Value = 5000
$ws.Range("M203").Value = 5000
$ws.Range("N203").Value = "$/docena de atados (3 kilos)"
$ws.Range("O203").Value = "Provincia de Cautín"
$ws.Range("P203").Value = 1667
$ws.Range("Q203").Value = 3
$ws.Range("R203").Value = "Hortaliza"

# Keep the date-formatted style used by the rest of column D.
$ws.Range("D203").NumberFormat = $ws.Range("D204").NumberFormat
